$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Default column width nudged very slightly by the resave (8.5703125 -> 8.578125) ---
$ws.StandardWidth = 8.578125

# --- Row 2: shared-string content changed (sample renamed, suffix renamed) ---
$ws.Range("A2").Value = "cerebellum_8rings"
$ws.Range("B2").Value = "064_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms"

# --- Rows 3-9: newly populated ring entries ---
$data = @(
    @{ Row = 3; Suffix = "064_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" },
    @{ Row = 4; Suffix = "065_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" },
    @{ Row = 5; Suffix = "066_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" },
    @{ Row = 6; Suffix = "067_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" },
    @{ Row = 7; Suffix = "068_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" },
    @{ Row = 8; Suffix = "069_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" },
    @{ Row = 9; Suffix = "070_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = "cerebellum_8rings"
    $ws.Cells.Item($r, 2).Value = $entry.Suffix
    $ws.Cells.Item($r, 3).Value = 1
    $ws.Cells.Item($r, 3).WrapText = $false
    $ws.Cells.Item($r, 4).Value = $r - 2
    $ws.Cells.Item($r, 5).Value = $r - 1
}

# --- Selection moves to B12 ---
$ws.Range("B12").Select()
